# Weekly update to the "Zapallo italiano" price dataset:
# two new price records are inserted at the top of the data block
# (rows 559-560), pushing the existing records (559-653) down by two
# rows (to 561-655). The row-fixed columns (A,B,C,E,F,G,H,I,R) are
# identical for every record in this block, so only the record-specific
# columns (D,J,K,L,M,N,O,P,Q) need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D=4 (Fecha), J=10 (Volumen), K=11 (Precio minimo), L=12 (Precio maximo),
# M=13 (Precio promedio ponderado), N=14 (Unidad de comercializacion),
# O=15 (Origen), P=16 (Precio $/Kg), Q=17 (Kg o Unidades)
$cols = @(4, 10, 11, 12, 13, 14, 15, 16, 17)

# A=1,B=2,C=3,E=5,F=6,G=7,H=8,I=9,R=18 are identical on every record of
# this block; rows 654-655 are brand new rows, so they need these too.
$fixedCols = @(1, 2, 3, 5, 6, 7, 8, 9, 18)

# Shift the existing 559..653 block down by two rows, working from the
# bottom up so we never overwrite a source row before it has been read.
for ($i = 655; $i -ge 561; $i--) {
    $src = $i - 2
    foreach ($c in $cols) {
        $ws.Cells.Item($i, $c).Value2 = $ws.Cells.Item($src, $c).Value2
    }
    if ($i -ge 654) {
        foreach ($c in $fixedCols) {
            $ws.Cells.Item($i, $c).Value2 = $ws.Cells.Item($src, $c).Value2
        }
        # Rows 654/655 are brand new cells - copy the date format used by
        # the rest of column D (style index 2, "YYYY-MM-DD HH:MM:SS").
        $ws.Cells.Item($i, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }
}

# Fill in the two brand-new records at rows 559 and 560.
$ws.Cells.Item(559, 4).Value2 = 45218
$ws.Cells.Item(559, 10).Value2 = 200
$ws.Cells.Item(559, 11).Value2 = 14000
$ws.Cells.Item(559, 12).Value2 = 14000
$ws.Cells.Item(559, 13).Value2 = 14000
$ws.Cells.Item(559, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(559, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(559, 16).Value2 = 280
$ws.Cells.Item(559, 17).Value2 = 50

$ws.Cells.Item(560, 4).Value2 = 45218
$ws.Cells.Item(560, 10).Value2 = 300
$ws.Cells.Item(560, 11).Value2 = 16000
$ws.Cells.Item(560, 12).Value2 = 16000
$ws.Cells.Item(560, 13).Value2 = 16000
$ws.Cells.Item(560, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(560, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(560, 16).Value2 = 320
$ws.Cells.Item(560, 17).Value2 = 50
